# Generate Report for Archive
# - Update localization status text from "Ready for handoff" to "In Translation"
#   on every sheet where it appears.
# - Shrink the column(s) that hold that status text now that the shorter
#   label no longer needs as much room.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
$newColumnWidth = 13.4101845877511

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count

    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $used.Cells.Item($r, $c)
            $cellValue = [string]$cell.Value()
            if ($cellValue -eq $oldStatus) {
                $cell.Value = $newStatus
                $ws.Columns.Item($cell.Column).ColumnWidth = $newColumnWidth
            }
        }
    }
}
